$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '41.089.83'
$ws.Range("E2").Value = '  +3.17%  '

$ws.Range("D3").Value = '2.242.93'
$ws.Range("E3").Value = '  +1.89%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '302.17'
$ws.Range("E5").Value = '  +3.24%  '

$ws.Range("D6").Value = '90.56'
$ws.Range("E6").Value = '  +4.53%  '

$ws.Range("D7").Value = '0.518'
$ws.Range("E7").Value = '  +2.10%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '0.481'
$ws.Range("E9").Value = '  +2.21%  '

$ws.Range("D10").Value = '53.99'
$ws.Range("E10").Value = '  +9.83%  '

$ws.Range("D11").Value = '31.67'
$ws.Range("E11").Value = '  +6.51%  '

$ws.Range("D12").Value = '0.0791'
$ws.Range("E12").Value = '  +2.27%  '

$ws.Range("E13").Value = '  +3.33%  '

$ws.Range("D14").Value = '6.54'
$ws.Range("E14").Value = '  +1.79%  '

$ws.Range("D15").Value = '2.592.25'
$ws.Range("E15").Value = '  +1.88%  '

$ws.Range("D16").Value = '14.03'
$ws.Range("E16").Value = '  +2.73%  '

$ws.Range("D17").Value = '2.265.11'
$ws.Range("E17").Value = '  +1.91%  '

$ws.Range("D18").Value = '0.747'
$ws.Range("E18").Value = '  +3.44%  '

$ws.Range("D19").Value = '41.032.16'
$ws.Range("E19").Value = '  +3.23%  '

$ws.Range("D20").Value = '11.80'
$ws.Range("E20").Value = '  +4.66%  '

$ws.Range("D21").Value = '0.0₃0899'
$ws.Range("E21").Value = '  +2.02%  '

$ws.Range("E22").Value = '  +1.96%  '

$ws.Range("D23").Value = '66.63'
$ws.Range("E23").Value = '  +2.42%  '

$ws.Range("D24").Value = '240.15'
$ws.Range("E24").Value = '  +1.97%  '

$ws.Range("E25").Value = '  +4.66%  '

$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("E27").Value = '  +2.47%  '

$ws.Range("D28").Value = '23.63'
$ws.Range("E28").Value = '  +5.92%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = '9.55'
$ws.Range("E29").Value = '  +4.58%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '2.09'
$ws.Range("E30").Value = '  -5.86%  '

$ws.Range("D31").Value = '157.58'
$ws.Range("E31").Value = '  +1.35%  '

$ws.Range("D32").Value = '32.98'
$ws.Range("E32").Value = '  +4.57%  '

$ws.Range("E33").Value = '  +0.08%  '

$ws.Range("D34").Value = '5.14'
$ws.Range("E34").Value = '  +5.69%  '

$ws.Range("D35").Value = '0.0728'
$ws.Range("E35").Value = '  +2.86%  '

$ws.Range("D36").Value = '3.00'
$ws.Range("E36").Value = '  +6.59%  '

$ws.Range("E37").Value = '  +1.25%  '

$ws.Range("B38").Value = 'Stellar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D38").Value = '0.115'
$ws.Range("E38").Value = '  +3.04%  '

$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").Value = '16.40'
$ws.Range("E39").Value = '  +6.35%  '

$ws.Range("E40").Value = '  +5.42%  '

$ws.Range("D41").Value = '1.76'
$ws.Range("E41").Value = '  +6.31%  '

$ws.Range("D42").Value = '3.90'
$ws.Range("E42").Value = '  +5.00%  '

$ws.Range("D43").Value = '2.071.82'
$ws.Range("E43").Value = '  -2.05%  '

$ws.Range("D44").Value = '19.95'
$ws.Range("E44").Value = '  +12.96%  '

$ws.Range("E45").Value = '  +3.53%  '

$ws.Range("D46").Value = '10.15'
$ws.Range("E46").Value = '  +5.66%  '

$ws.Range("D47").Value = '2.92'
$ws.Range("E47").Value = '  +11.16%  '

$ws.Range("E48").Value = '  -3.98%  '

$ws.Range("D49").Value = '2.464.30'
$ws.Range("E49").Value = '  +1.88%  '

$ws.Range("D50").Value = '1.49'
$ws.Range("E50").Value = '  +2.96%  '

$ws.Range("E51").Value = '  +3.94%  '
